$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.859.75"
$ws.Range("E2").Value = "  -8.28%  "
$ws.Range("D3").Value = "1.390.57"
$ws.Range("E3").Value = "  -9.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3010"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06349"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9540"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.189"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.72%  "
$ws.Range("D15").Value = "1.401.15"
$ws.Range("E15").Value = "  -8.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009803"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05598"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -14.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -15.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.435"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.248"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.54%  "
$ws.Range("D25").Value = "19.887.02"
$ws.Range("E25").Value = "  -8.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.45%  "
$ws.Range("D29").Value = "1.556.58"
$ws.Range("E29").Value = "  -8.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "106.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.859"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -20.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.144"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7753"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07518"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.214"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05541"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1862"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01980"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.016"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -13.34%  "
$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.07%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5103"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4914"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.692"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.016"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.42%  "
